$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.550.35'
$ws.Range("E2").Value = '  +0.83%  '

$ws.Range("D3").Value = '2.018.14'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'263.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.46%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = "'55.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.35%  '

$ws.Range("E9").Value = '  +0.58%  '

$ws.Range("D10").Value = "'0.0772"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.02%  '

$ws.Range("E11").Value = '  -2.09%  '

$ws.Range("D12").Value = "'14.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.47%  '

$ws.Range("D13").Value = '2.312.16'
$ws.Range("E13").Value = '  +0.59%  '

$ws.Range("E14").Value = '  -5.17%  '

$ws.Range("D15").Value = "'20.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.48%  '

$ws.Range("D16").Value = "'5.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.19%  '

$ws.Range("D17").Value = '2.009.84'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '37.480.98'
$ws.Range("E18").Value = '  +0.83%  '

$ws.Range("D19").Value = "'69.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.36%  '

$ws.Range("E20").Value = '  -2.68%  '

$ws.Range("E21").Value = '  -0.81%  '

$ws.Range("D22").Value = "'228.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.94%  '

$ws.Range("D23").Value = "'2.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.70%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.40%  '

$ws.Range("D26").Value = "'164.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.49%  '

$ws.Range("D27").Value = "'8.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.33%  '

$ws.Range("D28").Value = "'19.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("E29").Value = '  -11.01%  '

$ws.Range("E30").Value = '  -0.43%  '

$ws.Range("E31").Value = '  -1.21%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = "'0.0651"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.27%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'4.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.74%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("E35").Value = '  +1.49%  '

$ws.Range("E36").Value = '  +1.19%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").Value = "'3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").Value = "'5.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.63%  '

$ws.Range("E40").Value = '  +4.81%  '

$ws.Range("E41").Value = '  +2.65%  '

$ws.Range("E42").Value = '  -4.52%  '

$ws.Range("D43").Value = "'0.0213"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.06%  '

$ws.Range("D44").Value = '1.395.16'
$ws.Range("E44").Value = '  +1.09%  '

$ws.Range("D45").Value = "'90.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").Value = "'15.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.61%  '

$ws.Range("E47").Value = '  -2.14%  '

$ws.Range("D48").Value = "'7.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.15%  '

$ws.Range("D49").Value = "'2.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.88%  '

$ws.Range("D50").Value = '2.204.07'
$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("E51").Value = '  -4.39%  '
